$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: update ride
$ws.Range("D6").Value = "up_NotificationRosterEmailDevice 10, 1, 0"
$ws.Range("E6").Value = "up_NotificationRosterEmailDevice 10, 2, 0"

# Row 7: roster add
$ws.Range("D7").Value = "up_NotificationRosterEmailDevice 10, 1, 1"
$ws.Range("E7").Value = "up_NotificationRosterEmailDevice 10, 2, 1"

# Row 8: roster update
$ws.Range("D8").Value = "up_NotificationRosterEmailDevice 10, 1, 1"
$ws.Range("E8").Value = "up_NotificationRosterEmailDevice 10, 2, 1"

# Row 9: cancel ride
$ws.Range("D9").Value = "up_NotificationRosterEmailDevice 10, 1, 0"
$ws.Range("E9").Value = "up_NotificationRosterEmailDevice 10, 2, 0"

# Row 10: delete ride
$ws.Range("D10").Value = "up_NotificationRosterEmailDevice 10, 1, 0"
$ws.Range("E10").Value = "up_NotificationRosterEmailDevice 10, 2, 0"

# Row 11: chat add
$ws.Range("D11").Value = "up_NotificationRosterEmailDevice 10, 1, 1"
$ws.Range("E11").Value = "up_NotificationRosterEmailDevice 10, 2, 1"

# Row 12: friend add
$ws.Range("D12").Value = "up_NotificationFriendEmailDevice 1, 1"
$ws.Range("E12").Value = "up_NotificationFriendEmailDevice 1, 0"

# Row 13: friend request
$ws.Range("D13").Value = "up_NotificationFriendEmailDevice 1, 1"
$ws.Range("E13").Value = "up_NotificationFriendEmailDevice 1, 0"

# Row 14: hub member request
$ws.Range("D14").Value = "up_NotificationHubAdminEmailDevice 1, 1"
$ws.Range("E14").Value = "up_NotificationHubAdminEmailDevice 1, 2"

# Row 15: hub member add
$ws.Range("D15").Value = "up_NotificationFriendEmailDevice 1, 1"
$ws.Range("E15").Value = "up_NotificationFriendEmailDevice 1, 0"

# Row 16: new row - ride invite
$ws.Range("A16").Value = "ride invite"
$ws.Range("B16").Value = "selection"
$ws.Range("C16").Value = "ride detail"

# Update selection to D16
$ws.Range("D16").Select()
